$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# "Elimina EC anteriores y se agregan nuevos, se modifica base de datos"
# The previous account-statement periods (rows 16-22) are replaced: the list of
# periods is reversed (newest period now on top) and the "Salario Basico"
# (column G) is updated uniformly for every period.

$periodos = @("1908", "1907", "1906", "1905", "1904", "1903", "1901")
$valorMora = @(20979, 33125, 33125, 33125, 33125, 33125, 31249)
$salarioBasico = 828116

for ($i = 0; $i -lt 7; $i++) {
    $row = 16 + $i
    $ws.Cells.Item($row, 5).Value = $periodos[$i]
    $ws.Cells.Item($row, 6).Value = $valorMora[$i]
    $ws.Cells.Item($row, 7).Value = $salarioBasico
}
